$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 22:27:45"
$wsZhCn.Range("E3").Value = "2016-03-22 22:27:45"
$wsZhCn.Range("H2").Value = "2016-03-22 22:28:09"
$wsZhCn.Range("H3").Value = "2016-03-22 22:28:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 22:27:49"
$wsDeDe.Range("E3").Value = "2016-03-22 22:27:49"
$wsDeDe.Range("H2").Value = "2016-03-22 22:28:15"
$wsDeDe.Range("H3").Value = "2016-03-22 22:28:15"
